# Tidsplan_xjobb.xlsx update
# - Update ActlHours (P15) from 21 to 25 (downstream SUM/formula cells
#   recalc automatically).
# - Scroll/selection change: topLeftCell F1 -> G1, active cell P15 -> Q15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the actual-hours value for row 15.
$ws.Range("P15").Value = 25

# Update the view state: which column is left-most visible, and which
# cell is selected/active.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7           # column G is the 7th column
$ws.Range("Q15").Select()

$excel.CalculateFullRebuild()
